$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 19:19"

# Update statistics for Estados Unidos (row 4)
$ws.Range("B4").Value = 7203628
$ws.Range("C4").Value = 18157
$ws.Range("D4").Value = 4442259
$ws.Range("E4").Value = 2553494
$ws.Range("G4").Value = 337
$ws.Range("H4").Value = 207875

# Update statistics for India (row 5)
$ws.Range("B5").Value = 5877154
$ws.Range("C5").Value = 61051
$ws.Range("D5").Value = 4812155
$ws.Range("E5").Value = 972116
$ws.Range("G5").Value = 566
$ws.Range("H5").Value = 92883

# Update statistics for Turquia (row 21)
$ws.Range("B21").Value = 311455
$ws.Range("C21").Value = 1665
$ws.Range("D21").Value = 273282
$ws.Range("E21").Value = 30315
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 7858

# Update statistics for Chequia (row 57)
$ws.Range("B57").Value = 60027
$ws.Range("C57").Value = 1653
$ws.Range("D57").Value = 28022
$ws.Range("E57").Value = 31427
$ws.Range("G57").Value = 11
$ws.Range("H57").Value = 578

# Rows 74 & 75: Serbia and Libano swap ranking order and get new data
# Row 74 now holds Libano's (updated) data, row 75 now holds Serbia's (updated) data
$ws.Range("A74").Value = "Libano"
$ws.Range("B74").Value = 33962
$ws.Range("C74").Value = 1143
$ws.Range("D74").Value = 14778
$ws.Range("E74").Value = 18851
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 333

$ws.Range("A75").Value = "Serbia"
$ws.Range("B75").Value = 33238
$ws.Range("C75").Value = 75
$ws.Range("D75").Value = 31536
$ws.Range("E75").Value = 956
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 746

# Update statistics for Monaco (row 190)
$ws.Range("E190").Value = 37
$ws.Range("G190").Value = 1
$ws.Range("H190").Value = 2

# Update statistics for Antigua y Barbuda (row 196)
$ws.Range("B196").Value = 98
$ws.Range("C196").Value = 1
$ws.Range("E196").Value = 3

# Rows 215 & 216: Montserrat and Islas Malvinas swap ranking order and get new data
# Row 215 now holds Islas Malvinas's (updated) data, row 216 now holds Montserrat's (updated) data
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
